# Permite alterar registro padrao de frequencia
# e permite configurar numero padrao de aulas seguidas.
#
# Adds a new "Registro dos alunos" column (J) to the attendance sheet,
# pre-fills the existing rows with the default attendance code (0),
# and documents the code meanings (0/1/2) in a comment on the header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header in J1, matching the look of the other header cells (e.g. I1)
$ws.Range("J1").Value = "Registro dos alunos"
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Default attendance record (0 = Presente) for the existing data rows
$ws.Range("J2:J6").Value = 0

# Explain the code meaning on the new header cell
$commentText = "0: Presente" + [char]10 + "1: Ausente" + [char]10 + "2: Não registrado."
$ws.Range("J1").AddComment($commentText)

# Bring the new column into view / match the authored selection state
$ws.Range("I2").Select()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
